$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells we touch keep a Text format so numeric-looking
# strings (e.g. "594.52") are not auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.026.10"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.512.61"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.52"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.41"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  +3.00%  "
$ws.Range("E9").Value = "  +6.45%  "
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.119.74"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.71"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.015.39"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.495.79"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.32"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.22"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "394.75"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.27"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("E25").Value = "  -4.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.15"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.95"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.37"
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("E34").Value = "  +3.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.20"
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.896"
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.90"
$ws.Range("E38").Value = "  +2.55%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.61"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.798.36"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.86"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0304"
$ws.Range("E46").Value = "  -3.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "340.65"
$ws.Range("E47").Value = "  -3.74%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.46"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.845"
$ws.Range("E51").Value = "  -1.86%  "

# Row 41/42 swap of Coin name and Link, plus new Price/Volume values
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.38"
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.07"
$ws.Range("E42").Value = "  +2.13%  "
